$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Test Data"
$ws.Range("A2").Value = "etaeeigaehr"
$ws.Range("A3").Value = "rtrtrydmms"
$ws.Range("A4").Value = "eEeOUuoaRHRGSa"
$ws.Range("A5").Value = "iers-ta*_ta!d"
$ws.Range("A6").Value = "eouaeiuo"
$ws.Range("A7").Value = "e  gtr  w q  ii z "
$ws.Range("A8").Value = "1w1f1eg53qe4o"
$ws.Range("A9").ClearContents()

$ws.Range("A9").Select()
